# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I10").Value = "sv"
$ws.Range("J10").Value = "Statement-opinion"
$ws.Range("I15").Value = "sd"
$ws.Range("J15").Value = "Statement-non-opinion"
$ws.Range("I25").Value = "sd"
$ws.Range("J25").Value = "Statement-non-opinion"
$ws.Range("I34").Value = "sv"
$ws.Range("J34").Value = "Statement-opinion"
$ws.Range("I35").Value = "sd"
$ws.Range("J35").Value = "Statement-non-opinion"
$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I60").Value = "aa"
$ws.Range("J60").Value = "Agree/Accept"
$ws.Range("I76").Value = "ba"
$ws.Range("J76").Value = "Appreciation"
$ws.Range("I82").Value = "sd"
$ws.Range("J82").Value = "Statement-non-opinion"
$ws.Range("I85").Value = "sd"
$ws.Range("J85").Value = "Statement-non-opinion"
$ws.Range("I86").Value = "ba"
$ws.Range("J86").Value = "Appreciation"
$ws.Range("I88").Value = "sv"
$ws.Range("J88").Value = "Statement-opinion"
$ws.Range("I100").Value = "sv"
$ws.Range("J100").Value = "Statement-opinion"
$ws.Range("I104").Value = "%"
$ws.Range("J104").Value = "Uninterpretable"
$ws.Range("I110").Value = "b"
$ws.Range("J110").Value = "Acknowledge (Backchannel)"
$ws.Range("I116").Value = "sv"
$ws.Range("J116").Value = "Statement-opinion"
$ws.Range("I132").Value = "aa"
$ws.Range("J132").Value = "Agree/Accept"
$ws.Range("I134").Value = "sv"
$ws.Range("J134").Value = "Statement-opinion"
$ws.Range("I138").Value = "aa"
$ws.Range("J138").Value = "Agree/Accept"
$ws.Range("I147").Value = "sd"
$ws.Range("J147").Value = "Statement-non-opinion"
$ws.Range("I149").Value = "sd"
$ws.Range("J149").Value = "Statement-non-opinion"
$ws.Range("I150").Value = "sd"
$ws.Range("J150").Value = "Statement-non-opinion"
$ws.Range("I156").Value = "sd"
$ws.Range("J156").Value = "Statement-non-opinion"
$ws.Range("I165").Value = "sd"
$ws.Range("J165").Value = "Statement-non-opinion"
$ws.Range("I176").Value = "b"
$ws.Range("J176").Value = "Acknowledge (Backchannel)"
$ws.Range("I184").Value = "sd"
$ws.Range("J184").Value = "Statement-non-opinion"
$ws.Range("I187").Value = "sv"
$ws.Range("J187").Value = "Statement-opinion"
$ws.Range("I188").Value = "qy"
$ws.Range("J188").Value = "Yes-No-Question"
$ws.Range("I191").Value = "aa"
$ws.Range("J191").Value = "Agree/Accept"
$ws.Range("I197").Value = "sd"
$ws.Range("J197").Value = "Statement-non-opinion"
$ws.Range("I203").Value = "sv"
$ws.Range("J203").Value = "Statement-opinion"
$ws.Range("I220").Value = "b"
$ws.Range("J220").Value = "Acknowledge (Backchannel)"
$ws.Range("I226").Value = "sd"
$ws.Range("J226").Value = "Statement-non-opinion"
$ws.Range("I239").Value = "sd"
$ws.Range("J239").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "ba"
$ws.Range("J243").Value = "Appreciation"
$ws.Range("I273").Value = "sd"
$ws.Range("J273").Value = "Statement-non-opinion"
$ws.Range("I274").Value = "aa"
$ws.Range("J274").Value = "Agree/Accept"
$ws.Range("I280").Value = "sd"
$ws.Range("J280").Value = "Statement-non-opinion"
